# Append three new logbook rows (45, 46, 47) to Sheet1, extending the
# used range from A1:O44 to A1:O47, per the "Fixing the status column
# of the GUI" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=45; A=44; B="2025-11-30 00:34:56"; C="A605-010525-CHK-Y12"; D=3988; E=3988; F=3890; G=1; H=93; I=0; J=4; K=1963; L=$false; M=97;  N=2.43; O=11.12 },
    @{ Row=46; A=45; B="2025-11-30 00:49:07"; C="A873-150925-CHK-Y06"; D=891;  E=891;  F=881;  G=0; H=9;  I=0; J=1; K=631;  L=$false; M=10;  N=1.12; O=3.08  },
    @{ Row=47; A=46; B="2025-11-30 00:49:24"; C="A605-010525-CHK-Y12"; D=3988; E=3988; F=3890; G=1; H=93; I=0; J=4; K=1963; L=$false; M=97;  N=2.43; O=11.21 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value  = $r.A
    $ws.Cells.Item($n, 2).Value  = $r.B
    $ws.Cells.Item($n, 3).Value  = $r.C
    $ws.Cells.Item($n, 4).Value  = $r.D
    $ws.Cells.Item($n, 5).Value  = $r.E
    $ws.Cells.Item($n, 6).Value  = $r.F
    $ws.Cells.Item($n, 7).Value  = $r.G
    $ws.Cells.Item($n, 8).Value  = $r.H
    $ws.Cells.Item($n, 9).Value  = $r.I
    $ws.Cells.Item($n, 10).Value = $r.J
    $ws.Cells.Item($n, 11).Value = $r.K
    $ws.Cells.Item($n, 12).Value = $r.L
    $ws.Cells.Item($n, 13).Value = $r.M
    $ws.Cells.Item($n, 14).Value = $r.N
    $ws.Cells.Item($n, 15).Value = $r.O
}
